$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "303.61"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "3.85%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "31.98"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "7.79%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.247"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.72%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07547"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "5.23%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.929"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "5.20%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.821"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "6.35%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.517"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "8.09%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9226"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.37%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1694"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "4.32%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07913"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "3.15%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08050"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "3.48%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03051"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "4.56%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09916"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "10.22%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001513"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-5.18%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04596"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.74%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006340"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "3.56%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.448"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.96%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.231"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.11%"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.44%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1335"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-2.34%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.481"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "10.89%"

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.72%"

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.97%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004462"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "4.97%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001398"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "20.23%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001784"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "6.21%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01695"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "2,485.84%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04487"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.15%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006929"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.91%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1354"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "6.09%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002077"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01380"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "4.51%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006176"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "6.14%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.7191"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-62.73%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.01299"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.54%"
